$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook row-data (Fecha/Calidad/Volumen/Precios/Unidad/Origen/etc.) was
# reshuffled between rows. Snapshot source rows first (values only, cols D,L,M,N,O,P,Q,R,S,T),
# then write them into their destination rows, to correctly realize the permutation
# even though some rows are both a source and a destination.
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$snapshot = @{}
$snapshot[2] = @{}
foreach ($c in $cols) { $snapshot[2][$c] = $ws.Range("${c}2").Value2 }
$snapshot[3] = @{}
foreach ($c in $cols) { $snapshot[3][$c] = $ws.Range("${c}3").Value2 }
$snapshot[6] = @{}
foreach ($c in $cols) { $snapshot[6][$c] = $ws.Range("${c}6").Value2 }
$snapshot[7] = @{}
foreach ($c in $cols) { $snapshot[7][$c] = $ws.Range("${c}7").Value2 }
$snapshot[8] = @{}
foreach ($c in $cols) { $snapshot[8][$c] = $ws.Range("${c}8").Value2 }
$snapshot[9] = @{}
foreach ($c in $cols) { $snapshot[9][$c] = $ws.Range("${c}9").Value2 }
$snapshot[10] = @{}
foreach ($c in $cols) { $snapshot[10][$c] = $ws.Range("${c}10").Value2 }
$snapshot[11] = @{}
foreach ($c in $cols) { $snapshot[11][$c] = $ws.Range("${c}11").Value2 }
$snapshot[13] = @{}
foreach ($c in $cols) { $snapshot[13][$c] = $ws.Range("${c}13").Value2 }
$snapshot[14] = @{}
foreach ($c in $cols) { $snapshot[14][$c] = $ws.Range("${c}14").Value2 }
$snapshot[15] = @{}
foreach ($c in $cols) { $snapshot[15][$c] = $ws.Range("${c}15").Value2 }
$snapshot[16] = @{}
foreach ($c in $cols) { $snapshot[16][$c] = $ws.Range("${c}16").Value2 }
$snapshot[17] = @{}
foreach ($c in $cols) { $snapshot[17][$c] = $ws.Range("${c}17").Value2 }
$snapshot[18] = @{}
foreach ($c in $cols) { $snapshot[18][$c] = $ws.Range("${c}18").Value2 }
$snapshot[19] = @{}
foreach ($c in $cols) { $snapshot[19][$c] = $ws.Range("${c}19").Value2 }
$snapshot[20] = @{}
foreach ($c in $cols) { $snapshot[20][$c] = $ws.Range("${c}20").Value2 }

$destMap = @{
    2 = 11
    3 = 9
    6 = 10
    7 = 14
    8 = 20
    9 = 16
    10 = 3
    11 = 2
    13 = 19
    14 = 17
    15 = 6
    16 = 7
    17 = 8
    18 = 13
    19 = 18
    20 = 15
}

foreach ($src in $destMap.Keys) {
    $dst = $destMap[$src]
    foreach ($c in $cols) {
        $ws.Range("${c}$dst").Value = $snapshot[$src][$c]
    }
}
